$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> column letter -> new value
# Only the cells that actually change per the diff are listed here.
$updates = @{
    3  = @{ F = 12; H = 12 }
    6  = @{ E = 5;  F = 1;  H = 1 }
    15 = @{ E = 144; F = 75; H = 75 }
    18 = @{ E = 101; F = 42; H = 42 }
    28 = @{ E = 13; F = 9;  H = 9 }
    29 = @{ E = 15; F = 8;  H = 8 }
    32 = @{ E = 17 }
    36 = @{ F = 37; H = 37 }
    37 = @{ F = 23; H = 23 }
    38 = @{ E = 63; F = 13; H = 13 }
    42 = @{ E = 29 }
    44 = @{ E = 24; F = 11; H = 11 }
    48 = @{ F = 14; H = 14 }
    49 = @{ E = 56; F = 32; H = 32 }
    59 = @{ E = 8;  F = 4;  H = 4 }
    61 = @{ F = 8;  H = 8 }
    64 = @{ F = 15; H = 15 }
    73 = @{ F = 10; H = 10 }
    78 = @{ F = 16; H = 16 }
    82 = @{ E = 15 }
    88 = @{ F = 10; H = 10 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $cellRef = "$colLetter$rowNum"
        $ws.Range($cellRef).Value = $cols[$colLetter]
    }
}

$wb.Save()
